# Append new MTBF incident rows (221-230) to the worksheet, extending the
# used range from A1:H220 to A1:H230, as captured by the source commit
# ("Fixing MTBF data is saved and loaded correctly").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Bloque, Incidencia, Fecha, Hora, Turno, HoraReparacion, TiempoReparacion, MTBF
# $null placeholders mean "leave the cell empty" (mirrors existing sparse rows
# in the sheet where the repair time / duration were not recorded).
$data = @(
    @(221, "WC47 NACP", "Etiquetadora",          "2024-06-11", "12:32:40", "Mañana", "12:32:43", "0:00:03", "-0.01 minutos"),
    @(222, "WC47 NACP", "Fallo atornillador",     "2024-06-11", "12:32:53", "Mañana", "12:32:55", "0:00:02", "0.05 minutos"),
    @(223, "WC47 NACP", "No atornilla tapa",      "2024-06-11", "12:33:15", "Mañana", "12:33:24", "0:00:09", "-0.00 minutos"),
    @(224, "WC47 NACP", "Fallo tornillo",         "2024-06-11", "12:37:21", "Mañana", "12:37:23", "0:00:02", "-0.01 minutos"),
    @(225, "WC47 NACP", "No atornilla tapa",      "2024-06-11", "12:37:28", "Mañana", "12:37:31", "0:00:03", "0.03 minutos"),
    @(226, "WC47 NACP", "Fallo en elevador",      "2024-06-11", "12:38:20", "Mañana", "12:38:22", "0:00:02", "-0.00 minutos"),
    @(227, "WC47 NACP", "Fallo en paletizador",   "2024-06-11", "12:38:24", "Mañana", "12:38:26", "0:00:02", "0.02 minutos"),
    @(228, "WC47 NACP", "Etiquetadora",           "2024-06-11", "12:38:59", "Mañana", $null,       $null,     "-0.00 minutos"),
    @(229, "WC47 NACP", "Etiquetadora",           "2024-06-11", "12:40:29", "Mañana", "12:40:36", "0:00:07", "-0.00 minutos"),
    @(230, "WC47 NACP", "No coge placa",          "2024-06-11", "12:40:38", "Mañana", "12:40:42", "0:00:04", "0.04 minutos")
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]   # A: Bloque
    $ws.Cells.Item($r, 2).Value = $row[2]   # B: Incidencia

    # C: Fecha - force text so Excel does not turn "2024-06-11" into a date serial
    $ws.Range("C" + $r).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[3]

    $ws.Cells.Item($r, 4).Value = $row[4]   # D: Hora
    $ws.Cells.Item($r, 5).Value = $row[5]   # E: Turno

    if ($row[6] -ne $null) {
        $ws.Cells.Item($r, 6).Value = $row[6]   # F: Hora de Reparacion
    }
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 7).Value = $row[7]   # G: Tiempo de Reparacion
    }

    $ws.Cells.Item($r, 8).Value = $row[8]   # H: MTBF
}
